$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.727.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.63%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.197.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.19%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'601.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'157.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.80%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.198.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.19%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +4.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.51%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.512"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.91%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'39.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.722.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.17%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'66.724.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.61%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'7.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.18%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.197.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.10%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.112"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.64%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'515.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.18%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'15.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.54%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.44%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'14.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.50%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'85.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.50%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +2.24%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +9.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +9.22%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +9.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'28.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.51%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.76%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'519.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +10.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'54.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0898"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.36%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.19%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.125"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +7.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.80%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.26%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0₃0690"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +13.27%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +7.04%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.874.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +6.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = "'  +2.12%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +9.64%  "
$ws.Range("E51").Style = "Normal"

Write-Output "Applied 91 changes"
